$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 253, shifting rows 253:266 down to 254:267.
$ws.Rows.Item(253).Insert()

# Copy the (now shifted) row 254 values into the new row 253, then update
# the date (column D) to the new record's date.
$srcValues = $ws.Range("A254:R254").Value2
$ws.Range("A253:R253").Value2 = $srcValues

$ws.Range("D253").Value2 = 45013
